$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Correct the marks in the "Total" row (and related "Marking" row) of the marksheet.
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 85
$ws.Range("E12").Value = "85/140"
